$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4013624899316142
$ws.Range("C2").Value = 0.0866989264809348
$ws.Range("D2").Value = 0.05685525389132806
$ws.Range("F2").Value = 1.293505937744186
$ws.Range("G2").Value = 0.002486907507366913
$ws.Range("I2").Value = 0.9256184846368072
$ws.Range("K2").Value = 0.471969531767769
# Row 3
$ws.Range("B3").Value = 0.365269467197237
$ws.Range("C3").Value = 0.07840369955697213
$ws.Range("D3").Value = 0.05638767675157297
$ws.Range("F3").Value = 1.288992706591984
$ws.Range("G3").Value = 0.002490351079537989
$ws.Range("I3").Value = 0.9294360497602554
$ws.Range("K3").Value = 0.4288580102792707
# Row 4
$ws.Range("B4").Value = 0.3432446107302951
$ws.Range("C4").Value = 0.07334664119406398
$ws.Range("D4").Value = 0.0560976547413965
$ws.Range("F4").Value = 1.287022783986934
$ws.Range("G4").Value = 0.002492576795849499
$ws.Range("I4").Value = 0.9322897607712264
$ws.Range("K4").Value = 0.4025568541258906
# Row 5
$ws.Range("B5").Value = 0.3343038558852243
$ws.Range("C5").Value = 0.07129491054685388
$ws.Range("D5").Value = 0.05597875194122537
$ws.Range("F5").Value = 1.286421047345854
$ws.Range("G5").Value = 0.002493511880149283
$ws.Range("I5").Value = 0.9335805782021822
$ws.Range("K5").Value = 0.3918817059417279
# Row 6
$ws.Range("B6").Value = 0.3328213448391182
$ws.Range("C6").Value = 0.07095476803851852
$ws.Range("D6").Value = 0.05595896553586144
$ws.Range("F6").Value = 1.286333257534672
$ws.Range("G6").Value = 0.002493668849231351
$ws.Range("I6").Value = 0.9338026349901014
$ws.Range("K6").Value = 0.3901116949719494
# Row 7
$ws.Range("B7").Value = 0.3431238923405999
$ws.Range("C7").Value = 0.0733189342135887
$ws.Range("D7").Value = 0.05609605405140883
$ws.Range("F7").Value = 1.28701385542508
$ws.Range("G7").Value = 0.002492589293127402
$ws.Range("I7").Value = 0.9323066516503502
$ws.Range("K7").Value = 0.4024127118777869
# Row 8
$ws.Range("B8").Value = 0.3888894432214443
$ws.Range("C8").Value = 0.08383117537914586
$ws.Range("D8").Value = 0.05669465251120087
$ws.Range("F8").Value = 1.291783207072939
$ws.Range("G8").Value = 0.002488071793053551
$ws.Range("I8").Value = 0.9268288310016715
$ws.Range("K8").Value = 0.4570695852218876
# Row 9
$ws.Range("B9").Value = 0.4797117918444656
$ws.Range("C9").Value = 0.1047367203502176
$ws.Range("D9").Value = 0.05784447508269608
$ws.Range("F9").Value = 1.307515661617316
$ws.Range("G9").Value = 0.002480092513742202
$ws.Range("I9").Value = 0.920143615526996
$ws.Range("K9").Value = 0.5655951296696458
# Row 10
$ws.Range("B10").Value = 0.5470943528563339
$ws.Range("C10").Value = 0.1202799263139696
$ws.Range("D10").Value = 0.0586736337612308
$ws.Range("F10").Value = 1.322997742861219
$ws.Range("G10").Value = 0.002474760729564521
$ws.Range("I10").Value = 0.917723132723161
$ws.Range("K10").Value = 0.6461551734052762
# Row 11
$ws.Range("B11").Value = 0.5778911211501736
$ws.Range("C11").Value = 0.1273923636536267
$ws.Range("D11").Value = 0.05904725697877211
$ws.Range("F11").Value = 1.330900415883931
$ws.Range("G11").Value = 0.002472449184818105
$ws.Range("I11").Value = 0.9171668508520341
$ws.Range("K11").Value = 0.6829856340172
# Row 12
$ws.Range("B12").Value = 0.5895736524929021
$ws.Range("C12").Value = 0.1300917495734382
$ws.Range("D12").Value = 0.05918820876921416
$ws.Range("F12").Value = 1.334017137197293
$ws.Range("G12").Value = 0.002471590153005635
$ws.Range("I12").Value = 0.9170348498118841
$ws.Range("K12").Value = 0.6969587192143081
# Row 13
$ws.Range("B13").Value = 0.5870567042807124
$ws.Range("C13").Value = 0.1295101180863014
$ws.Range("D13").Value = 0.05915787619490231
$ws.Range("F13").Value = 1.333340365708665
$ws.Range("G13").Value = 0.002471774437060808
$ws.Range("I13").Value = 0.9170597756768331
$ws.Range("K13").Value = 0.6939482031623072
# Row 14
$ws.Range("B14").Value = 0.5788518413769168
$ws.Range("C14").Value = 0.1276143219782284
$ws.Range("D14").Value = 0.05905886391764881
$ws.Range("F14").Value = 1.331154339325735
$ws.Range("G14").Value = 0.002472378185428415
$ws.Range("I14").Value = 0.9171544130037219
$ws.Range("K14").Value = 0.6841346842223857
# Row 15
$ws.Range("B15").Value = 0.5738287833140419
$ws.Range("C15").Value = 0.12645388302127
$ws.Range("D15").Value = 0.05899814631855804
$ws.Range("F15").Value = 1.32983151831894
$ws.Range("G15").Value = 0.002472750119216219
$ws.Range("I15").Value = 0.9172226328973494
$ws.Range("K15").Value = 0.67812702620202
# Row 16
$ws.Range("B16").Value = 0.545084589157284
$ws.Range("C16").Value = 0.119815956786681
$ws.Range("D16").Value = 0.05864914326769366
$ws.Range("F16").Value = 1.322498629757234
$ws.Range("G16").Value = 0.002474914079623144
$ws.Range("I16").Value = 0.9177704763258916
$ws.Range("K16").Value = 0.6437518900102646
# Row 17
$ws.Range("B17").Value = 0.5274876800709478
$ws.Range("C17").Value = 0.1157545380352758
$ws.Range("D17").Value = 0.05843411471276383
$ws.Range("F17").Value = 1.318220719843339
$ws.Range("G17").Value = 0.00247627071746501
$ws.Range("I17").Value = 0.9182463182481087
$ws.Range("K17").Value = 0.622710681719667
# Row 18
$ws.Range("B18").Value = 0.5173799822134413
$ws.Range("C18").Value = 0.1134224441723859
$ws.Range("D18").Value = 0.05831010147544546
$ws.Range("F18").Value = 1.315841076714264
$ws.Range("G18").Value = 0.002477061746151591
$ws.Range("I18").Value = 0.9185712714236374
$ws.Range("K18").Value = 0.6106255973346038
# Row 19
$ws.Range("B19").Value = 0.5139600295179036
$ws.Range("C19").Value = 0.1126335099975222
$ws.Range("D19").Value = 0.05826805584674588
$ws.Range("F19").Value = 1.315049249030508
$ws.Range("G19").Value = 0.002477331419800906
$ws.Range("I19").Value = 0.9186900898068018
$ws.Range("K19").Value = 0.6065367654026659
# Row 20
$ws.Range("B20").Value = 0.5293594976137967
$ws.Range("C20").Value = 0.1161864761859022
$ws.Range("D20").Value = 0.05845703962500437
$ws.Range("F20").Value = 1.318667734345922
$ws.Range("G20").Value = 0.002476125191541497
$ws.Range("I20").Value = 0.9181903564542182
$ws.Range("K20").Value = 0.6249487689360365
# Row 21
$ws.Range("B21").Value = 0.581261255160257
$ws.Range("C21").Value = 0.1281709983263397
$ws.Range("D21").Value = 0.05908796078660572
$ws.Range("F21").Value = 1.331793054770571
$ws.Range("G21").Value = 0.002472200408134716
$ws.Range("I21").Value = 0.9171244787115427
$ws.Range("K21").Value = 0.6870164431686874
# Row 22
$ws.Range("B22").Value = 0.6153012149040933
$ws.Range("C22").Value = 0.1360389267421169
$ws.Range("D22").Value = 0.05949719966268674
$ws.Range("F22").Value = 1.341095092313736
$ws.Range("G22").Value = 0.002469730307889923
$ws.Range("I22").Value = 0.9168864460031045
$ws.Range("K22").Value = 0.7277338319668445
# Row 23
$ws.Range("B23").Value = 0.5971226286385161
$ws.Range("C23").Value = 0.1318364120014337
$ws.Range("D23").Value = 0.05927907123960097
$ws.Range("F23").Value = 1.336064015323004
$ws.Range("G23").Value = 0.002471039983302604
$ws.Range("I23").Value = 0.9169714267192859
$ws.Range("K23").Value = 0.7059882888015068
# Row 24
$ws.Range("B24").Value = 0.5285132208525454
$ws.Range("C24").Value = 0.1159911880144193
$ws.Range("D24").Value = 0.05844667648667468
$ws.Range("F24").Value = 1.318465390597922
$ws.Range("G24").Value = 0.002476190949205148
$ws.Range("I24").Value = 0.9182154967354919
$ws.Range("K24").Value = 0.6239368929830675
# Row 25
$ws.Range("B25").Value = 0.4550268145255529
$ws.Range("C25").Value = 0.09904939102375465
$ws.Range("D25").Value = 0.05753609188373332
$ws.Range("F25").Value = 1.302573082676474
$ws.Range("G25").Value = 0.002482157541004066
$ws.Range("I25").Value = 0.9215159732039808
$ws.Range("K25").Value = 0.5360915260727097
